$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$notes = $s.NotesPage
$shape = $notes.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = ""
